# Adapt column header formatting to respective input file names
# Rename "<Name>_old" -> "<Name>_FV2310" and "<Name>_new" -> "<Name>_FV2404"
# then expose the header row as an Excel Table (ListObject) and freeze the
# header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$baseNames = @("Segmentname","Segmentgruppe","Segment","Datenelement","Segment ID","Code","Qualifier","Beschreibung","Bedingungsausdruck","Bedingung")

# Columns A-J (1-10) carried the "_old" suffix -> becomes "_FV2310"
for ($i = 0; $i -lt $baseNames.Length; $i++) {
    $col = $i + 1
    $ws.Cells.Item(1, $col).Value = $baseNames[$i] + "_FV2310"
}

# Column K (11) is "diff" - untouched

# Columns L-U (12-21) carried the "_new" suffix -> becomes "_FV2404"
for ($i = 0; $i -lt $baseNames.Length; $i++) {
    $col = $i + 12
    $ws.Cells.Item(1, $col).Value = $baseNames[$i] + "_FV2404"
}

# Turn the used range into a proper Excel Table ("Table1")
$tbl = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $ws.Range("A1:U70"), $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$tbl.Name = "Table1"
$tbl.TableStyle = ""

# Freeze the header row (row 1)
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
